# Add 2022-Q4 data
#
# 1) Insert a brand-new worksheet named "2022-Q4" right after the "总计"
#    (summary) sheet, pushing every existing quarter sheet down by one
#    position. Populate it with the new quarter's fund-holding data.
# 2) Update the "总计" summary sheet: insert a new row right under the
#    header for "2022-Q4" (count=1, value=0.08) and push the previously
#    existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q4" worksheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# Header row (bold, centered, thin box border - matches the sibling
# quarter sheets).
$header = $q4.Range("B1:H1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Data row - same box-border/bold styling is applied to column A like
# the other quarter sheets.
$aCell = $q4.Range("A2")
$aCell.Font.Bold = $true
$aCell.HorizontalAlignment = -4108
$aCell.VerticalAlignment = -4160
$aCell.Borders.LineStyle = 1
$aCell.Value = 0

# These columns hold numeric-looking values that must stay TEXT (matches
# the source data format used on every other quarter sheet) - a leading
# quote forces text without leaving behind a Text number-format override.
$q4.Range("B2").Value = "'561550"
$q4.Range("C2").Value = "华泰柏瑞中证500增强策略ETF"
$q4.Range("D2").Value = "'7.54"
$q4.Range("E2").Value = "'99.26"
$q4.Range("F2").Value = "'1.12"
$q4.Range("G2").Value = "'0.0844"
$q4.Range("H2").Value = 6

# Restore the active/selected tab to the last sheet ("2020-Q4"), which
# is where it was before this edit.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()

# ---------------------------------------------------------------------
# 2. "总计" summary sheet - add the 2022-Q4 row
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()

# Inserting a row copies the formatting of the row above (the bold
# header) onto B2:D2 - drop that so the new data row is plain, like
# every other data row in this sheet. Column A keeps the header-row
# style, so copy that format over explicitly from the row below.
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy($total.Range("A2"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.08

# Leave the originally-active sheet selected.
$lastSheet.Activate()
